$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 554: revised daily figures (SI / intubated / hors-SI patients) ---
$ws.Range("E554").Value = 10
$ws.Range("F554").Value = 6
$ws.Range("G554").Value = 18

# --- Row 555: revised intubated count ---
$ws.Range("F555").Value = 6

# --- Row 561: revised new-cases count ---
$ws.Range("C561").Value = 119

# --- Row 562: revised new-cases + hors-SI count ---
$ws.Range("C562").Value = 78
$ws.Range("G562").Value = 14

# --- Row 563: day now in the past -> fill in the real observed data ---
$ws.Range("C563").Value = 74
$ws.Range("E563").Value = 11
$ws.Range("F563").Value = 9
$ws.Range("G563").Value = 12
$ws.Range("L563").Value = 0
$ws.Range("M563").Value = 0

# --- Row 564 ---
$ws.Range("C564").Value = 67
$ws.Range("E564").Value = 11
$ws.Range("F564").Value = 8
$ws.Range("G564").Value = 14
$ws.Range("L564").Value = 0
$ws.Range("M564").Value = 0

# --- Row 565 ---
$ws.Range("C565").Value = 40
$ws.Range("E565").Value = 11
$ws.Range("F565").Value = 9
$ws.Range("G565").Value = 15
$ws.Range("L565").Value = 0
$ws.Range("M565").Value = 0

# --- Row 566 ---
$ws.Range("C566").Value = 4
$ws.Range("E566").Value = 11
$ws.Range("F566").Value = 8
$ws.Range("G566").Value = 16
$ws.Range("L566").Value = 0
$ws.Range("M566").Value = 0

# --- Restore the view: scrolled back to the top of the frozen pane, cursor on A2 ---
$ws.Range("A2").Select()
